# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1568
    4  = 1021
    5  = 22
    7  = 2610
    8  = 39
    9  = 1646
    11 = 67
    12 = 538
    15 = 58
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
